$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "Arthritis Care & Research"
$ws.Range("G3").Value = "https://openalex.org/S13698240"
$ws.Range("H3").Value = "Wiley"
$ws.Range("I3").Value = "2151-464X"
